$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the Area / Atotal columns (G, H)
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Row 2: base area segment (relative to 0) and running total
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Row 3: first incremental area segment (its own formula, not yet shared)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share one formula pattern (fill-down), anchored relative to G4
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

$ws.Range("F2").Select() | Out-Null
